# Auto-generated edit script: updates the cryptos price/volume table
# on Sheet1 to match the "Updated cryptos list" GitHub Actions commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '30.809.33'
$ws.Range('E2').Value = '  +1.36%  '

# Row 3
$ws.Range('D3').Value = '1.886.84'
$ws.Range('E3').Value = '  +2.07%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9985'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.21%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '239.10'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +2.46%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9987'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.20%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4776'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +2.21%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2873'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +5.26%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06572'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +4.44%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.94'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +16.10%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '97.53'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +16.30%  '

# Row 12
$ws.Range('D12').Value = '1.872.70'
$ws.Range('E12').Value = '  +1.15%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07576'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +1.48%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.130'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +4.04%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6560'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +5.74%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '309.78'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +35.45%  '

# Row 17
$ws.Range('D17').Value = '30.799.26'
$ws.Range('E17').Value = '  +1.54%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.20'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +6.64%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.9984'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.21%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007587'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +3.84%  '

# Row 21
$ws.Range('D21').Value = '2.122.72'
$ws.Range('E21').Value = '  +2.58%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9990'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.14%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.129'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +4.42%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.194'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +5.49%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.319'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +1.88%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '166.68'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.16%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.34'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +14.10%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.950'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +4.22%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1074'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +5.14%  '

# Row 30
$ws.Range('E30').Value = '  -1.72%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.174'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +2.24%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.982'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +4.35%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05042'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +3.45%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.178'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +3.42%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7369'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +5.25%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.709'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.71%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01959'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +1.97%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.709'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +1.83%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.082'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +7.41%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9048'
$ws.Range('D40').ClearFormats()

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '107.96'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +2.01%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9983'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.21%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4227'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +4.92%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.645'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +2.34%  '

# Row 45
$ws.Range('E45').Value = '  +7.65%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.384'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +4.34%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.1228'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +1.94%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.011'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +4.53%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '34.81'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +4.38%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05617'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +1.38%  '

# Row 51
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.391'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +3.43%  '

